$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Commit: cell C10 value changed from 18 to 1 (restored from an earlier revision).
$ws.Range("C10").Value = 1
